# Expands the Berserker AI spec: splits the title/"Berserker" mentions so
# Word's spell-checker brackets them with proofErr marks, inserts the new
# state-machine / fuzzy-logic / perception paragraphs, and relocates the
# _GoBack bookmark to the end of the (now much longer) perception paragraph.
#
# The safest way to get the exact run layout (including <w:proofErr/> spell
# -check markers and the <w:bookmarkStart/End> pair) is to hand Word a
# WordOpenXML fragment for the whole story range instead of poking at
# Find/Replace run-by-run.

$d = $word.ActiveDocument

$newBodyXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Necesidades del </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:t>Berserker</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t xml:space="preserve"> según sus mecánicas:</w:t></w:r></w:p><w:p><w:r><w:t>-Máquina de estados: Cambiar entre estado Caminar,</w:t></w:r><w:r><w:t xml:space="preserve"> Investigar,</w:t></w:r><w:r><w:t xml:space="preserve"> Luchar, Huir</w:t></w:r></w:p><w:p><w:r><w:t>Estado inicial: Caminar. Deambula por una cierta zona de la nave.</w:t></w:r></w:p><w:p><w:r><w:t>Cosas que puede hacer en este estado: escuchar, investigar.</w:t></w:r></w:p><w:p><w:r><w:t>Pasa a estado investigar si escucha un ruido de intensidad media o mayor dentro de su rango de e</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t>cucha.</w:t></w:r></w:p><w:p><w:r><w:t>Estado investigar: Al oír un ruido raro, irá a su origen. Si ve al jugador, alertará a sus aliados y pasará al estado atacar. Si no ve al jugador, vuelve a donde estaba y sigue caminando.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Estado atacar: Después de </w:t></w:r><w:r><w:t xml:space="preserve">ver al jugador, el </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:t>Berserker</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t xml:space="preserve"> empezará a luchar, prestando atención a la distancia con el jugador y su vida.</w:t></w:r></w:p><w:p><w:r><w:t>Cosas que puede hacer en este estado: Correr, ataque uñas, ataque ácido.</w:t></w:r></w:p><w:p><w:r><w:t>Pasa a estado huir si su vida baja al 30%.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Estado huir: El </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:t>Berserker</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t xml:space="preserve"> decide abandonar la pelea, y huye del jugador. Usará cualquier ruta disponible para alejarse corriendo del jugador. Si no puede huir, se queda ahí hasta que el jugador lo mate. Si consigue alejarse una cierta distancia del jugador, se quedará ahí.</w:t></w:r></w:p><w:p><w:r><w:t>-Lógica difusa: decidir la distancia a la que puede ver u oír. Decidir el momento en el que tiene poca vida y debe huir</w:t></w:r><w:r><w:t xml:space="preserve"> (30% de su vida máxima)</w:t></w:r><w:r><w:t>. Decidir el ataque a realizar en una lucha según la distancia</w:t></w:r><w:r><w:t>(Cerca=menos de 1 metro. Medio=Entre 1 y 6 metros. Lejos=Más de 6 metros)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>-Percepción sensorial: vista y oído.</w:t></w:r><w:r><w:t xml:space="preserve"> La vista es un arco de 30 grados delante del </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:t>Ber</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t>erker</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t xml:space="preserve">, con una distancia máxima de 10 metros. No puede ver a través de las paredes. El oído es un circulo alrededor del </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:t>Berserker</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Su radio para oír sonidos medios es 15 metros. Su radio para </w:t></w:r><w:r><w:t>oír</w:t></w:r><w:r><w:t xml:space="preserve"> sonidos altos es 30 metros. No puede oír sonidos bajos. Hay que tener en cuenta que la intensidad de un sonido se puede reducir un nivel si hay paredes o puertas cerradas de por medio (es decir, el </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:t>berserker</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t xml:space="preserve"> no puede oír un sonido medio al otro lado de una </w:t></w:r><w:r><w:t>pared,</w:t></w:r><w:r><w:t xml:space="preserve"> aunque esté dentro de su rango de oído</w:t></w:r><w:r><w:t>)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack" /><w:bookmarkEnd w:id="0" /><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>-</w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:t>Pathfinding</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:t>Dijkstra</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t>: para saber por dónde puede andar, para saber huir del jugador.</w:t></w:r></w:p><w:p /><w:p /><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Replace the whole main story (everything except the trailing sectPr,
# which Document.Content never includes) with the fully-specified XML.
$d.Content.InsertXML($newBodyXml)
